$d = $word.ActiveDocument

function Replace-ParagraphXml {
    param(
        [string]$FindText,
        [string]$ParaAttrs,
        [string]$InnerXml
    )

    $rng = $d.Content
    $ok = $rng.Find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find text: $FindText"
    }

    $para = $rng.Paragraphs(1).Range

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body><w:p ' + $ParaAttrs + '>' + $InnerXml + '</w:p></w:body>' + `
        '</w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'

    $para.InsertXML($xml)
}


$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="en-US"/></w:rPr>'

# --- Edit 1: "Defines the kind of data a variable can hold. They can be further " ---
$inner1 = `
    '<w:pPr>' + $rPr + '</w:pPr>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">Defines the kind of data a </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">variable </w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>holds</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">. They can be further </w:t></w:r>' + `
    '<w:r w:rsidR="000708DF">' + $rPr + '<w:t>classifi</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>ed into two c</w:t></w:r>' + `
    '<w:r w:rsidR="000708DF">' + $rPr + '<w:t>ategories:</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> primitive and non-primitive. Primitive data types</w:t></w:r>' + `
    '<w:r w:rsidR="000708DF">' + $rPr + '<w:t xml:space="preserve"> are immutable and represent a single value. They are stored by value (stored directly). Non-primitive data types</w:t></w:r>' + `
    '<w:r w:rsidR="00CD7A6B">' + $rPr + '<w:t xml:space="preserve"> are objects that are stored by reference. </w:t></w:r>'

$attrs1 = 'w14:paraId="5D307471" w14:textId="596BAAD0" w:rsidR="00E074C2" w:rsidRDefault="00E074C2" w:rsidP="00E074C2"'
Replace-ParagraphXml "Defines the kind of data a variable can hold. They can be further " $attrs1 $inner1

# --- Edit 2: "Refers to where the variable to accessible " ---
$inner2 = `
    '<w:pPr>' + $rPr + '</w:pPr>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">Refers to where the variable </w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>is</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> accessible </w:t></w:r>' + `
    '<w:r w:rsidR="00D52FFA">' + $rPr + '<w:t>in the code. It is essential for understanding code and avoiding bugs.</w:t></w:r>'

$attrs2 = 'w14:paraId="5EBD50F1" w14:textId="3C89B83E" w:rsidR="00327849" w:rsidRDefault="00327849" w:rsidP="00E074C2"'
Replace-ParagraphXml "Refers to where the variable to accessible " $attrs2 $inner2

Write-Host "Done applying edits"
